$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.132.80'
$ws.Range("E2").Value = '  +4.09%  '
$ws.Range("D3").Value = '2.340.07'
$ws.Range("E3").Value = '  +1.92%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.50'
$ws.Range("E5").Value = '  +2.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.88'
$ws.Range("E6").Value = '  +3.94%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.539'
$ws.Range("E8").Value = '  +1.59%  '
$ws.Range("D9").Value = '2.354.70'
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("E10").Value = '  +6.43%  '
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.24'
$ws.Range("E12").Value = '  +3.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.341'
$ws.Range("E13").Value = '  +0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.67'
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("D15").Value = '2.744.39'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").Value = '56.954.59'
$ws.Range("E16").Value = '  +3.74%  '
$ws.Range("E17").Value = '  +2.13%  '
$ws.Range("D18").Value = '2.343.90'
$ws.Range("E18").Value = '  +2.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.48'
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.02'
$ws.Range("E22").Value = '  -2.42%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.07'
$ws.Range("E24").Value = '  +0.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.163'
$ws.Range("E25").Value = '  +7.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.993'
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.85'
$ws.Range("E27").Value = '  +4.24%  '
$ws.Range("E28").Value = '  +10.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.76'
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("E30").Value = '  +5.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.71'
$ws.Range("E31").Value = '  +4.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.18'
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.30'
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +0.54%  '
$ws.Range("E36").Value = '  +1.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.920'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.02'
$ws.Range("E38").Value = '  +3.68%  '
$ws.Range("E39").Value = '  +7.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.92'
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("E42").Value = '  +4.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.92'
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '278.17'
$ws.Range("E44").Value = '  +7.92%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.15'
$ws.Range("E45").Value = '  +0.77%  '
$ws.Range("E46").Value = '  +2.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0504'
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("E49").Value = '  +4.41%  '
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.52'
$ws.Range("E51").Value = '  +5.64%  '
